$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.172.75"
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "2.055.23"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.12%  "

$ws.Range("E6").Value = "  -1.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.31"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.06%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0782"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.109"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.75%  "

$ws.Range("D13").Value = "2.356.27"
$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("E14").Value = "  +1.75%  "

$ws.Range("E15").Value = "  +2.61%  "

$ws.Range("D16").Value = "2.057.62"
$ws.Range("E16").Value = "  -0.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +15.71%  "

$ws.Range("D18").Value = "37.177.51"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").Value = "0.0₃0899"
$ws.Range("E20").Value = "  -3.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "237.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  +1.80%  "

$ws.Range("E25").Value = "  -9.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.83%  "

$ws.Range("E29").Value = "  -1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0618"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.55%  "

$ws.Range("E33").Value = "  +0.80%  "

$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("E37").Value = "  -0.77%  "

$ws.Range("E38").Value = "  -2.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.43%  "

$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +16.03%  "

$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.103"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.98%  "

$ws.Range("E42").Value = "  -2.36%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.49%  "

$ws.Range("E46").Value = "  -1.53%  "

$ws.Range("E47").Value = "  -1.53%  "

$ws.Range("D48").Value = "1.275.96"
$ws.Range("E48").Value = "  -2.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.15%  "

$ws.Range("D50").Value = "2.240.62"
$ws.Range("E50").Value = "  -0.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.96%  "
